$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.979488
$ws.Range("H2").Value = 8.938464
$ws.Range("I2").Value = 0.01229331913219231
$ws.Range("J2").Value = 0.01233795916068688
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.02264433333333333
$ws.Range("N2").Value = 0.06793299999999999
$ws.Range("O2").Value = 0.08454793804489194
$ws.Range("P2").Value = 0.08454793804489193
$ws.Range("Q2").Value = 0.06746851943466667
$ws.Range("R2").Value = 0.6072166749119999
$ws.Range("S2").Value = 0.00103937478435468
$ws.Range("T2").Value = 0.001043149006718161

$ws.Range("G3").Value = 2.979488
$ws.Range("H3").Value = 8.938464
$ws.Range("I3").Value = 0.01229331913219231
$ws.Range("J3").Value = 0.01233795916068688
$ws.Range("M3").Value = 0.245184
$ws.Range("N3").Value = 0.735552
$ws.Range("O3").Value = 0.915452061955108
$ws.Range("P3").Value = 0.9154520619551081
$ws.Range("Q3").Value = 0.7305227857919999
$ws.Range("R3").Value = 6.574705072127999
$ws.Range("S3").Value = 0.01125394434783763
$ws.Range("T3").Value = 0.01129481015396872

$ws.Range("H4").Value = 42.55606899999999
$ws.Range("I4").Value = 0.05852855000910628
$ws.Range("J4").Value = 0.05874108139400382
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02264433333333333
$ws.Range("N4").Value = 0.06793299999999999
$ws.Range("O4").Value = 0.08454793804489194
$ws.Range("P4").Value = 0.08454793804489193
$ws.Range("Q4").Value = 0.3212179372641111
$ws.Range("R4").Value = 2.890961435376999
$ws.Range("S4").Value = 0.004948468220027277
$ws.Range("T4").Value = 0.004966437310390189

$ws.Range("H5").Value = 42.55606899999999
$ws.Range("I5").Value = 0.05852855000910628
$ws.Range("J5").Value = 0.05874108139400382
$ws.Range("M5").Value = 0.245184
$ws.Range("N5").Value = 0.735552
$ws.Range("O5").Value = 0.915452061955108
$ws.Range("P5").Value = 0.9154520619551081
$ws.Range("Q5").Value = 3.478022407231999
$ws.Range("R5").Value = 31.302201665088
$ws.Range("S5").Value = 0.053580081789079
$ws.Range("T5").Value = 0.05377464408361363

$ws.Range("G6").Value = 109.4796983333333
$ws.Range("H6").Value = 328.439095
$ws.Range("I6").Value = 0.4517114585149561
$ws.Range("J6").Value = 0.4533517325664633
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.02264433333333333
$ws.Range("N6").Value = 0.06793299999999999
$ws.Range("O6").Value = 0.08454793804489194
$ws.Range("P6").Value = 0.08454793804489193
$ws.Range("Q6").Value = 2.479094782292778
$ws.Range("R6").Value = 22.311853040635
$ws.Range("S6").Value = 0.03819127240869028
$ws.Range("T6").Value = 0.03832995419757375

$ws.Range("G7").Value = 109.4796983333333
$ws.Range("H7").Value = 328.439095
$ws.Range("I7").Value = 0.4517114585149561
$ws.Range("J7").Value = 0.4533517325664633
$ws.Range("M7").Value = 0.245184
$ws.Range("N7").Value = 0.735552
$ws.Range("O7").Value = 0.915452061955108
$ws.Range("P7").Value = 0.9154520619551081
$ws.Range("Q7").Value = 26.84267035616
$ws.Range("R7").Value = 241.58403320544
$ws.Range("S7").Value = 0.4135201861062658
$ws.Range("T7").Value = 0.4150217783688895

$ws.Range("G8").Value = 2.6307215
$ws.Range("H8").Value = 5.261443
$ws.Range("I8").Value = 0.01085431421352248
$ws.Range("J8").Value = 0.007262485910362437
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.02264433333333333
$ws.Range("N8").Value = 0.06793299999999999
$ws.Range("O8").Value = 0.08454793804489194
$ws.Range("P8").Value = 0.08454793804489193
$ws.Range("Q8").Value = 0.05957093455316666
$ws.Range("R8").Value = 0.357425607319
$ws.Range("S8").Value = 0.0009177098856446884
$ws.Range("T8").Value = 0.0006140282088012239

$ws.Range("G9").Value = 2.6307215
$ws.Range("H9").Value = 5.261443
$ws.Range("I9").Value = 0.01085431421352248
$ws.Range("J9").Value = 0.007262485910362437
$ws.Range("M9").Value = 0.245184
$ws.Range("N9").Value = 0.735552
$ws.Range("O9").Value = 0.915452061955108
$ws.Range("P9").Value = 0.9154520619551081
$ws.Range("Q9").Value = 0.645010820256
$ws.Range("R9").Value = 3.870064921536
$ws.Range("S9").Value = 0.009936604327877789
$ws.Range("T9").Value = 0.006648457701561214

$ws.Range("G10").Value = 113.0911763333333
$ws.Range("H10").Value = 339.273529
$ws.Range("I10").Value = 0.4666123581302228
$ws.Range("J10").Value = 0.4683067409684837
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.02264433333333333
$ws.Range("N10").Value = 0.06793299999999999
$ws.Range("O10").Value = 0.08454793804489194
$ws.Range("P10").Value = 0.08454793804489193
$ws.Range("Q10").Value = 2.560874293950778
$ws.Range("R10").Value = 23.047868645557
$ws.Range("S10").Value = 0.03945111274617501
$ws.Range("T10").Value = 0.03959436932140861

$ws.Range("G11").Value = 113.0911763333333
$ws.Range("H11").Value = 339.273529
$ws.Range("I11").Value = 0.4666123581302228
$ws.Range("J11").Value = 0.4683067409684837
$ws.Range("M11").Value = 0.245184
$ws.Range("N11").Value = 0.735552
$ws.Range("O11").Value = 0.915452061955108
$ws.Range("P11").Value = 0.9154520619551081
$ws.Range("Q11").Value = 27.728146978112
$ws.Range("R11").Value = 249.553322803008
$ws.Range("S11").Value = 0.4271612453840478
$ws.Range("T11").Value = 0.4287123716470751
